$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weight Data")

function Set-TextCell($cell, $text) {
    # Force a literal-text value (matches how the source data already
    # stores its Date column - as text, not as a real date serial) while
    # avoiding leaving a date-ish number format behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# New weekly weigh-in rows appended after the existing data (last row was 99).
$rows = @(
    @{ Row = 100; Date = "2025-10-13"; Weight = 197;   Name = "Marat";   Fat = $null;  Waist = $null },
    @{ Row = 101; Date = "2025-10-13"; Weight = 207.1; Name = "Michael"; Fat = 14.94;  Waist = 39 },
    @{ Row = 102; Date = "2025-10-20"; Weight = 199.8; Name = "Marat";   Fat = $null;  Waist = $null },
    @{ Row = 103; Date = "2025-10-27"; Weight = 209.5; Name = "Michael"; Fat = 15.38;  Waist = 39 },
    @{ Row = 104; Date = "2025-10-28"; Weight = 201.2; Name = "Marat";   Fat = $null;  Waist = $null }
)

foreach ($r in $rows) {
    # Column A (Date) holds text that looks like a date ("yyyy-mm-dd"), so it
    # needs the text-forcing helper or Excel auto-converts it to a real date.
    Set-TextCell $ws.Cells.Item($r.Row, 1) $r.Date
    $ws.Cells.Item($r.Row, 2).Value = $r.Weight
    # Column C (Name) is plain text - a normal Value assignment is enough.
    $ws.Cells.Item($r.Row, 3).Value = $r.Name
    if ($r.Fat -ne $null) {
        $ws.Cells.Item($r.Row, 4).Value = $r.Fat
    }
    if ($r.Waist -ne $null) {
        $ws.Cells.Item($r.Row, 5).Value = $r.Waist
    }
}
